# Update "想去人数" (column F) counts on both the "展览" and "全部类型"
# worksheets, which carry duplicate data in this workbook.

$wb = $excel.ActiveWorkbook

# Row number -> new value for column F
$updates = @{
    4  = 123
    5  = 114
    6  = 477
    16 = 18
    20 = 1026
    21 = 1429
    23 = 344
    30 = 268
    34 = 103
    38 = 3814
    39 = 4
    40 = 444
    42 = 957
    43 = 63
    46 = 84
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
